$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(1, 1).Value = -0.46475278436945189
$ws.Cells.Item(1, 2).Value = 0.46297903389995554
$ws.Cells.Item(2, 1).Value = -0.36837178422067396
$ws.Cells.Item(2, 2).Value = 0.36218180723016324
$ws.Cells.Item(3, 1).Value = -0.1855141216786933
$ws.Cells.Item(3, 2).Value = 0.18426093660362852
$ws.Cells.Item(4, 1).Value = -0.17226093661854058
$ws.Cells.Item(4, 2).Value = 0.17111543970966281
$ws.Cells.Item(5, 1).Value = -0.16511543975825038
$ws.Cells.Item(5, 2).Value = 0.16280370070739902
$ws.Cells.Item(6, 1).Value = -0.061764083062814557
$ws.Cells.Item(6, 2).Value = 0.061708958800599678
$ws.Cells.Item(7, 1).Value = -0.041708958861400802
$ws.Cells.Item(7, 2).Value = 0.041617136087170437
$ws.Cells.Item(8, 1).Value = -0.039506505505944922
$ws.Cells.Item(8, 2).Value = 0.039376602456405507
$ws.Cells.Item(9, 1).Value = -0.03337660250776775
$ws.Cells.Item(9, 2).Value = 0.033281023422827793
$ws.Cells.Item(10, 1).Value = -0.02728102347460748
$ws.Cells.Item(10, 2).Value = 0.02726927030363413
$ws.Cells.Item(11, 1).Value = -0.022769270354295656
$ws.Cells.Item(11, 2).Value = 0.022753742370305474
$ws.Cells.Item(12, 1).Value = -0.016753742422195295
$ws.Cells.Item(12, 2).Value = 0.016721505914028878
$ws.Cells.Item(13, 1).Value = -0.010721505966242439
$ws.Cells.Item(13, 2).Value = 0.010716978363936036
$ws.Cells.Item(14, 1).Value = 0.0012830215793782074
$ws.Cells.Item(14, 2).Value = -0.0012843991699877932
$ws.Cells.Item(15, 1).Value = -0.02105549061710299
$ws.Cells.Item(15, 2).Value = 0.021028862658519998
$ws.Cells.Item(16, 1).Value = -0.015028862710960045
$ws.Cells.Item(16, 2).Value = 0.015004979411355368
$ws.Cells.Item(17, 1).Value = -0.0090049794640503222
$ws.Cells.Item(17, 2).Value = 0.0089999999450345314
$ws.Cells.Item(18, 1).Value = -0.040956007396975025
$ws.Cells.Item(18, 2).Value = 0.040935087900656697
$ws.Cells.Item(19, 1).Value = -0.03193508794789679
$ws.Cells.Item(19, 2).Value = 0.031815163028243276
$ws.Cells.Item(20, 1).Value = -0.018014393137491425
$ws.Cells.Item(20, 2).Value = 0.018004370528887037
$ws.Cells.Item(21, 1).Value = -0.0090043705767053339
$ws.Cells.Item(21, 2).Value = 0.0089999999521177543
$ws.Cells.Item(22, 1).Value = -0.16766768560077239
$ws.Cells.Item(22, 2).Value = 0.16648952509683923
$ws.Cells.Item(23, 1).Value = -0.08465009951342406
$ws.Cells.Item(23, 2).Value = 0.084130122745670555
$ws.Cells.Item(24, 1).Value = -0.042130122821006033
$ws.Cells.Item(24, 2).Value = 0.041999999924217768
$ws.Cells.Item(25, 1).Value = -0.095039617694943956
$ws.Cells.Item(25, 2).Value = 0.094783157622686076
$ws.Cells.Item(26, 1).Value = -0.088783157673223201
$ws.Cells.Item(26, 2).Value = 0.088454609816967178
$ws.Cells.Item(27, 1).Value = -0.082454609867846251
$ws.Cells.Item(27, 2).Value = 0.081338114890909363
$ws.Cells.Item(28, 1).Value = -0.07533811494292042
$ws.Cells.Item(28, 2).Value = 0.074564244196997365
$ws.Cells.Item(29, 1).Value = -0.062564244254270207
$ws.Cells.Item(29, 2).Value = 0.062179543403157567
$ws.Cells.Item(30, 1).Value = -0.042179543466736469
$ws.Cells.Item(30, 2).Value = 0.042022130249439904
$ws.Cells.Item(31, 1).Value = -0.027022130309935832
$ws.Cells.Item(31, 2).Value = 0.027001216172456211
$ws.Cells.Item(32, 1).Value = -0.0060012162374940559
$ws.Cells.Item(32, 2).Value = 0.0059999999459776632
